# Update countries & provincias Spain
# Applies the data refresh: new stats for Australia / Finlandia, new case
# data for Fiyi, the corresponding re-sort of several country rows that
# shifted position, and the refreshed "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Footer timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 2 de Abril de 2020 a las 06:20"

# Australia (row 22): new confirmed/active/death totals
$ws.Cells.Item(22, 2).Value = 5106
$ws.Cells.Item(22, 3).Value = 58
$ws.Cells.Item(22, 5).Value = 4738

# Finlandia (row 42): active/recovered counts updated
$ws.Cells.Item(42, 4).Value = 300
$ws.Cells.Item(42, 5).Value = 1129

# Countries re-sorted around rows 185-197 (data follows each country as
# it moves to its new row position); row 185 "Fiyi" also gained new cases.
$ws.Cells.Item(185, 1).Value = "Fiyi"
$ws.Cells.Item(185, 3).Value = 2

$ws.Cells.Item(186, 1).Value = "Republica del Chad"
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 7
$ws.Cells.Item(186, 8).Value = 0

$ws.Cells.Item(187, 1).Value = "Sudan"
$ws.Cells.Item(187, 2).Value = 7
$ws.Cells.Item(187, 4).Value = 2
$ws.Cells.Item(187, 5).Value = 3
$ws.Cells.Item(187, 8).Value = 2

$ws.Cells.Item(188, 1).Value = "Santa Sede"

$ws.Cells.Item(189, 1).Value = "Liberia"

$ws.Cells.Item(190, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(190, 5).Value = 6
$ws.Cells.Item(190, 8).Value = 0

$ws.Cells.Item(192, 1).Value = "Cabo Verde"
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 5).Value = 5

$ws.Cells.Item(193, 1).Value = "Mauritania"
$ws.Cells.Item(193, 2).Value = 6
$ws.Cells.Item(193, 4).Value = 2
$ws.Cells.Item(193, 5).Value = 3
$ws.Cells.Item(193, 8).Value = 1

$ws.Cells.Item(194, 1).Value = "Nepal"
$ws.Cells.Item(194, 4).Value = 1
$ws.Cells.Item(194, 8).Value = 0

$ws.Cells.Item(197, 1).Value = "Nicaragua"
$ws.Cells.Item(197, 4).Value = 0
$ws.Cells.Item(197, 8).Value = 1

# Two more country swaps further down the table (stats identical, only
# the name moves to its new row).
$ws.Cells.Item(200, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(201, 1).Value = "Belice"

$ws.Cells.Item(203, 1).Value = "Burundi"
$ws.Cells.Item(204, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(206, 1).Value = "Anguila"
